$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K to F:M)
$ws.Columns("D:E").Insert()

# Copy number formatting/style from column F (the old column D, now shifted) into new D and E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new columns with the latest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 140000
$ws.Range("E8").Value = 166400
$ws.Range("D9").Value = 80100
$ws.Range("E9").Value = 99200
$ws.Range("D10").Value = 59900
$ws.Range("E10").Value = 67200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 107700
$ws.Range("E17").Value = 127300
$ws.Range("D18").Value = 32300
$ws.Range("E18").Value = 39100
$ws.Range("D20").Value = 500
$ws.Range("E20").Value = 200
$ws.Range("D21").Value = 36200
$ws.Range("E21").Value = 43100
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("D23").Value = 32800
$ws.Range("E23").Value = 39300
$ws.Range("D24").Value = 7600
$ws.Range("E24").Value = 9800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 25200
$ws.Range("E26").Value = 29500
$ws.Range("D27").Value = 25200
$ws.Range("E27").Value = 29500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("E32").Value = -200
$ws.Range("D33").Value = 25200
$ws.Range("E33").Value = 29500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 25200
$ws.Range("E35").Value = 29500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 105700
$ws.Range("E41").Value = 107300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 99600
$ws.Range("E43").Value = 96900
$ws.Range("D44").Value = 57800
$ws.Range("E44").Value = 35500
$ws.Range("D45").Value = 7100
$ws.Range("E45").Value = 9100
$ws.Range("D46").Value = 270200
$ws.Range("E46").Value = 248800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 117100
$ws.Range("E48").Value = 108200
$ws.Range("D49").Value = 74500
$ws.Range("E49").Value = 74600
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3300
$ws.Range("E52").Value = 3300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 465100
$ws.Range("E54").Value = 434900
$ws.Range("D57").Value = 31100
$ws.Range("E57").Value = 18200
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 61700
$ws.Range("E59").Value = 57600
$ws.Range("D60").Value = 92800
$ws.Range("E60").Value = 75800
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 29400
$ws.Range("E62").Value = 30400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 122200
$ws.Range("E66").Value = 106300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 416900
$ws.Range("E72").Value = 391800
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 343000
$ws.Range("E76").Value = 328700
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 25200
$ws.Range("E81").Value = 29500
$ws.Range("D83").Value = 3400
$ws.Range("E83").Value = 3800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 22800
$ws.Range("E89").Value = 117000
$ws.Range("D91").Value = -12200
$ws.Range("E91").Value = -3900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -12200
$ws.Range("E94").Value = -3900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -12200
$ws.Range("E100").Value = -8700
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = -1600
$ws.Range("E102").Value = 104400
